$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Goal-adjustment fixes: D32 now derives from C3 instead of a hard literal.
$ws.Range("D32").Formula = "=5-C3"

# Totals in the first table now also roll up row 32 (previously stopped at 31).
$ws.Range("C33").Formula = "=SUM(C28:C32)"
$ws.Range("D33").Formula = "=SUM(D28:D32)"

# Updated weight for row 41 in the second table.
$ws.Range("E41").Value = 1

# New rollup row 47 for the F/G columns, mirroring the existing C47/D47 pattern.
$ws.Range("F47").Formula = "=-6.5+F44"
$ws.Range("G47").Formula = "=-6+G44"

# Move the active selection to where the author left off.
[void]$ws.Range("D43").Select()
